$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new day column (07-nov) before the
#     existing "01-oct." column (DJ), shifting DJ:EN -> DK:EO. ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("DJ1").EntireColumn.Insert()
$ws1.Range("DJ1").Value = "07-nov"
$ws1.Range("DJ2:DJ25").Value = "-"

# --- Sheet "Gaz": append row 143 for 2025-11-05. ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A143").Value = "'2025-11-05"
$ws2.Range("A143").Style = "Normal"
$ws2.Range("B143").Value = 30.425

# --- Sheet "CO2": append row 143 for 2025-11-05. ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A143").Value = "'2025-11-05"
$ws3.Range("A143").Style = "Normal"
$ws3.Range("B143").Value = 81.18000000000001
